function Set-TextValue {
    param($ws, $cellRef, $value)
    $range = $ws.Range($cellRef)
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws "D2" "27.459.65"
Set-TextValue $ws "E2" "  -2.40%  "
Set-TextValue $ws "D3" "1.740.61"
Set-TextValue $ws "E3" "  -3.34%  "
Set-TextValue $ws "E4" "  -0.32%  "
Set-TextValue $ws "D5" "323.80"
Set-TextValue $ws "E5" "  -4.21%  "
Set-TextValue $ws "D6" "1.000"
Set-TextValue $ws "E6" "  -0.20%  "
Set-TextValue $ws "D7" "0.4251"
Set-TextValue $ws "E7" "  -8.74%  "
Set-TextValue $ws "E8" "  -3.63%  "
Set-TextValue $ws "E9" "  -1.21%  "
Set-TextValue $ws "D10" "1.126"
Set-TextValue $ws "E10" "  -1.81%  "
Set-TextValue $ws "D11" "0.07429"
Set-TextValue $ws "E11" "  -3.43%  "
Set-TextValue $ws "D12" "1.001"
Set-TextValue $ws "E12" "  -0.37%  "
Set-TextValue $ws "D13" "21.71"
Set-TextValue $ws "E13" "  -3.36%  "
Set-TextValue $ws "D14" "6.066"
Set-TextValue $ws "E14" "  -5.60%  "
Set-TextValue $ws "D15" "7.174"
Set-TextValue $ws "E15" "  -2.97%  "
Set-TextValue $ws "D16" "1.736.95"
Set-TextValue $ws "E16" "  -3.46%  "
Set-TextValue $ws "D17" "0.00001065"
Set-TextValue $ws "E17" "  -2.95%  "
Set-TextValue $ws "D18" "87.12"
Set-TextValue $ws "E18" "  +5.96%  "
Set-TextValue $ws "D19" "0.05972"
Set-TextValue $ws "E19" "  -11.78%  "
Set-TextValue $ws "D20" "0.9995"
Set-TextValue $ws "E20" "  -0.28%  "
Set-TextValue $ws "D21" "16.89"
Set-TextValue $ws "E21" "  -3.56%  "
Set-TextValue $ws "D22" "6.067"
Set-TextValue $ws "E22" "  -6.07%  "
Set-TextValue $ws "D23" "0.5230"
Set-TextValue $ws "E23" "  -5.21%  "
Set-TextValue $ws "D24" "27.483.67"
Set-TextValue $ws "E24" "  -2.37%  "
Set-TextValue $ws "E25" "  -5.74%  "
Set-TextValue $ws "D26" "2.407"
Set-TextValue $ws "E26" "  +0.13%  "
Set-TextValue $ws "D27" "20.15"
Set-TextValue $ws "E27" "  -3.59%  "
Set-TextValue $ws "D28" "2.370"
Set-TextValue $ws "E28" "  -1.75%  "
Set-TextValue $ws "D29" "150.16"
Set-TextValue $ws "E29" "  -0.85%  "
Set-TextValue $ws "D30" "1.937.16"
Set-TextValue $ws "E30" "  -3.68%  "
Set-TextValue $ws "D31" "1.269"
Set-TextValue $ws "E31" "  +0.07%  "
Set-TextValue $ws "E32" "  -6.26%  "
Set-TextValue $ws "D33" "3.720"
Set-TextValue $ws "E33" "  -8.37%  "
Set-TextValue $ws "D34" "5.604"
Set-TextValue $ws "E34" "  -6.36%  "
Set-TextValue $ws "E35" "  -6.94%  "
Set-TextValue $ws "D36" "12.56"
Set-TextValue $ws "E36" "  +2.88%  "
Set-TextValue $ws "D37" "0.2158"
Set-TextValue $ws "E37" "  -3.30%  "
Set-TextValue $ws "B38" "Hedera"
Set-TextValue $ws "C38" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws "D38" "0.06152"
Set-TextValue $ws "E38" "  -3.59%  "
Set-TextValue $ws "B39" "VeChain"
Set-TextValue $ws "C39" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws "D39" "0.02278"
Set-TextValue $ws "E39" "  -4.51%  "
Set-TextValue $ws "D40" "0.6426"
Set-TextValue $ws "E40" "  -4.27%  "
Set-TextValue $ws "D41" "5.024"
Set-TextValue $ws "E41" "  -4.55%  "
Set-TextValue $ws "D42" "1.184"
Set-TextValue $ws "E42" "  -3.70%  "
Set-TextValue $ws "D43" "1.423"
Set-TextValue $ws "E43" "  -4.42%  "
Set-TextValue $ws "D44" "0.9997"
Set-TextValue $ws "E44" "  -0.27%  "
Set-TextValue $ws "D45" "7.838"
Set-TextValue $ws "E45" "  -3.30%  "
Set-TextValue $ws "D46" "13.55"
Set-TextValue $ws "E46" "  -5.05%  "
Set-TextValue $ws "D47" "3.746"
Set-TextValue $ws "E47" "  -3.20%  "
Set-TextValue $ws "D48" "0.5857"
Set-TextValue $ws "E48" "  -5.11%  "
Set-TextValue $ws "D49" "125.44"
Set-TextValue $ws "E49" "  -3.30%  "
Set-TextValue $ws "D50" "1.941"
Set-TextValue $ws "E50" "  -5.71%  "
Set-TextValue $ws "D51" "0.06827"
Set-TextValue $ws "E51" "  -4.11%  "
